$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Remove the old "_GoBack" bookmark from its current location
#    (the empty paragraph before "Руководитель").
# ---------------------------------------------------------------------
$oldBm = $d.Bookmarks.Item("_GoBack")
$oldBm.Delete()

# ---------------------------------------------------------------------
# 2. Locate the final paragraph (the one with the date and the
#    reviewer's name) and replace the trailing "<tab>Name" run with a
#    tab run + a separate text run bearing the new name, preserving the
#    <w:tab/> element rather than collapsing it into literal text.
# ---------------------------------------------------------------------
$lastParaIndex = $d.Paragraphs.Count
$p = $d.Paragraphs.Item($lastParaIndex)
$pEnd = $p.Range.End - 1          # position right before the paragraph mark

$oldName = "Гринченко Н.Н."
$newName = "Чичикин В.А."
$runStart = $pEnd - $oldName.Length - 1   # -1 to include the preceding tab
$runRange = $d.Range($runStart, $pEnd)

$xmlFragment = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + `
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
  '<pkg:xmlData>' + `
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
  '<w:body><w:p>' + `
  '<w:r><w:tab/></w:r>' + `
  '<w:r><w:t>' + $newName + '</w:t></w:r>' + `
  '</w:p></w:body></w:document>' + `
  '</pkg:xmlData></pkg:part></pkg:package>'

$runRange.InsertXML($xmlFragment)

# ---------------------------------------------------------------------
# 3. Re-insert the "_GoBack" bookmark at the very end of that last
#    paragraph (after the new name, before the paragraph mark). A
#    placeholder character is used because collapsed Bookmarks.Add
#    calls exactly at a paragraph-end boundary are not honoured.
# ---------------------------------------------------------------------
$p = $d.Paragraphs.Item($lastParaIndex)
$endPos = $p.Range.End - 1
$ip = $d.Range($endPos, $endPos)
$ip.InsertAfter("X")

$p = $d.Paragraphs.Item($lastParaIndex)
$bmPos = ($p.Range.End - 1) - 1
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

$p = $d.Paragraphs.Item($lastParaIndex)
$delPos = $p.Range.End - 1
$delRange = $d.Range($delPos - 1, $delPos)
$delRange.Delete()
